$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5803
$ws.Range("E2").Value = 983
$ws.Range("F2").Value = 983
$ws.Range("G2").Value = 543
$ws.Range("H2").Value = 719
$ws.Range("I2").Value = 717
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 13007
$ws.Range("L2").Value = 8018
$ws.Range("M2").Value = 4989
$ws.Range("N2").Value = 4932
$ws.Range("O2").Value = 57
$ws.Range("P2").Value = 1221
$ws.Range("Q2").Value = 1750
$ws.Range("R2").Value = -1068
$ws.Range("S2").Value = -894
$ws.Range("T2").Value = 1099
$ws.Range("U2").Value = 652
$ws.Range("V2").Value = 7035
$ws.Range("W2").Value = 16.94
$ws.Range("X2").Value = 12.39
$ws.Range("AA2").Value = 160.72
$ws.Range("AB2").Value = 89.01000000000001
$ws.Range("AC2").Value = 2958
$ws.Range("AD2").Value = 7.91
$ws.Range("AE2").Value = 20199
$ws.Range("AF2").Value = 1.16
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 24416108
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()

# Row 3
$ws.Range("D3").Value = 5317
$ws.Range("E3").Value = 860
$ws.Range("F3").Value = 860
$ws.Range("G3").Value = 404
$ws.Range("H3").Value = 395
$ws.Range("I3").Value = 389
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 15179
$ws.Range("L3").Value = 9486
$ws.Range("M3").Value = 5694
$ws.Range("N3").Value = 5632
$ws.Range("O3").Value = 62
$ws.Range("P3").Value = 1221
$ws.Range("Q3").Value = 1383
$ws.Range("R3").Value = -2513
$ws.Range("S3").Value = 1138
$ws.Range("T3").Value = 2872
$ws.Range("U3").Value = -1490
$ws.Range("V3").Value = 8821
$ws.Range("W3").Value = 16.18
$ws.Range("X3").Value = 7.43
$ws.Range("Y3").Value = 7.37
$ws.Range("Z3").Value = 2.8
$ws.Range("AA3").Value = 166.6
$ws.Range("AB3").Value = 120.8
$ws.Range("AC3").Value = 1594
$ws.Range("AD3").Value = 11.92
$ws.Range("AE3").Value = 23056
$ws.Range("AF3").Value = 0.82
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 24426504

# Row 4
$ws.Range("D4").Value = 5403
$ws.Range("E4").Value = 441
$ws.Range("F4").Value = 441
$ws.Range("G4").Value = 351
$ws.Range("H4").Value = 308
$ws.Range("I4").Value = 298
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 22045
$ws.Range("L4").Value = 15784
$ws.Range("M4").Value = 6260
$ws.Range("N4").Value = 6066
$ws.Range("O4").Value = 194
$ws.Range("P4").Value = 1221
$ws.Range("Q4").Value = 1207
$ws.Range("R4").Value = -3743
$ws.Range("S4").Value = 2773
$ws.Range("T4").Value = 4287
$ws.Range("U4").Value = -3080
$ws.Range("V4").Value = 13930
$ws.Range("W4").Value = 8.16
$ws.Range("X4").Value = 5.7
$ws.Range("Y4").Value = 5.1
$ws.Range("Z4").Value = 1.65
$ws.Range("AA4").Value = 252.14
$ws.Range("AB4").Value = 144.66
$ws.Range("AC4").Value = 1220
$ws.Range("AD4").Value = 14.1
$ws.Range("AE4").Value = 24834
$ws.Range("AF4").Value = 0.6899999999999999
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 24426504

# Row 5
$ws.Range("D5").Value = 15607
$ws.Range("E5").Value = 1009
$ws.Range("F5").Value = 1009
$ws.Range("G5").Value = 1164
$ws.Range("H5").Value = 1092
$ws.Range("I5").Value = 1267
$ws.Range("J5").Value = -175
$ws.Range("K5").Value = 25868
$ws.Range("L5").Value = 18990
$ws.Range("M5").Value = 6878
$ws.Range("N5").Value = 6621
$ws.Range("O5").Value = 258
$ws.Range("P5").Value = 1221
$ws.Range("Q5").Value = 2273
$ws.Range("R5").Value = -7767
$ws.Range("S5").Value = 5231
$ws.Range("T5").Value = 7574
$ws.Range("U5").Value = -5301
$ws.Range("V5").Value = 16939
$ws.Range("W5").Value = 6.46
$ws.Range("X5").Value = 7
$ws.Range("Y5").Value = 19.97
$ws.Range("Z5").Value = 4.56
$ws.Range("AA5").Value = 276.08
$ws.Range("AB5").Value = 248.85
$ws.Range("AC5").Value = 5186
$ws.Range("AD5").Value = 4.41
$ws.Range("AE5").Value = 27104
$ws.Range("AF5").Value = 0.84
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 24426504

# Row 6
$ws.Range("D6").Value = 13347
$ws.Range("E6").Value = 1437
$ws.Range("F6").Value = 1437
$ws.Range("G6").Value = 841
$ws.Range("H6").Value = 870
$ws.Range("I6").Value = 852
$ws.Range("K6").Value = 29559
$ws.Range("L6").Value = 21540
$ws.Range("M6").Value = 8019
$ws.Range("N6").Value = 7720
$ws.Range("P6").Value = 1221
$ws.Range("Q6").Value = 2607
$ws.Range("R6").Value = -3088
$ws.Range("S6").Value = 791
$ws.Range("T6").Value = 2780
$ws.Range("U6").Value = -173
$ws.Range("V6").Value = 18701
$ws.Range("W6").Value = 10.77
$ws.Range("X6").Value = 6.52
$ws.Range("Y6").Value = 11.88
$ws.Range("Z6").Value = 3.14
$ws.Range("AA6").Value = 268.62
$ws.Range("AB6").Value = 320.66
$ws.Range("AC6").Value = 3488
$ws.Range("AD6").Value = 6.24
$ws.Range("AE6").Value = 31895
$ws.Range("AF6").Value = 0.68
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 24426504
$ws.Range("AI6").ClearContents()

# Row 7
$ws.Range("D7").Value = 10833
$ws.Range("E7").Value = 1612
$ws.Range("G7").Value = 1083
$ws.Range("H7").Value = 1136
$ws.Range("I7").Value = 1087
$ws.Range("K7").Value = 33294
$ws.Range("L7").Value = 23811
$ws.Range("M7").Value = 9485
$ws.Range("N7").Value = 9058
$ws.Range("P7").Value = 1220
$ws.Range("Q7").Value = 2889
$ws.Range("R7").Value = -3536
$ws.Range("S7").Value = 723
$ws.Range("T7").Value = 3262
$ws.Range("U7").Value = 76
$ws.Range("W7").Value = 14.88
$ws.Range("X7").Value = 10.49
$ws.Range("Y7").Value = 12.96
$ws.Range("Z7").Value = 3.62
$ws.Range("AA7").Value = 251.04
$ws.Range("AC7").Value = 4451
$ws.Range("AD7").Value = 4.63
$ws.Range("AE7").Value = 37760
$ws.Range("AF7").Value = 0.55
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 13010
$ws.Range("E8").Value = 1782
$ws.Range("G8").Value = 1158
$ws.Range("H8").Value = 1062
$ws.Range("I8").Value = 1028
$ws.Range("K8").Value = 35774
$ws.Range("L8").Value = 25128
$ws.Range("M8").Value = 10648
$ws.Range("N8").Value = 10076
$ws.Range("P8").Value = 1220
$ws.Range("Q8").Value = 2985
$ws.Range("R8").Value = -3268
$ws.Range("S8").Value = -55
$ws.Range("T8").Value = 3147
$ws.Range("U8").Value = 655
$ws.Range("W8").Value = 13.69
$ws.Range("X8").Value = 8.16
$ws.Range("Y8").Value = 10.74
$ws.Range("Z8").Value = 3.07
$ws.Range("AA8").Value = 235.99
$ws.Range("AC8").Value = 4207
$ws.Range("AD8").Value = 4.9
$ws.Range("AE8").Value = 42000
$ws.Range("AF8").Value = 0.49
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 13529
$ws.Range("E9").Value = 1860
$ws.Range("G9").Value = 1212
$ws.Range("H9").Value = 1107
$ws.Range("I9").Value = 1072
$ws.Range("K9").Value = 37896
$ws.Range("L9").Value = 26039
$ws.Range("M9").Value = 11858
$ws.Range("N9").Value = 11102
$ws.Range("P9").Value = 1220
$ws.Range("Q9").Value = 3198
$ws.Range("R9").Value = -3324
$ws.Range("S9").Value = -50
$ws.Range("T9").Value = 3280
$ws.Range("U9").Value = 615
$ws.Range("W9").Value = 13.75
$ws.Range("X9").Value = 8.18
$ws.Range("Y9").Value = 10.13
$ws.Range("Z9").Value = 3
$ws.Range("AA9").Value = 219.6
$ws.Range("AC9").Value = 4390
$ws.Range("AD9").Value = 4.69
$ws.Range("AE9").Value = 46278
$ws.Range("AF9").Value = 0.45
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").ClearContents()
